$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 21
$ws.Range("E2").Value = 30.07999992370605
$ws.Range("F2").Value = 31.45999908447266
$ws.Range("G2").Value = 21
$ws.Range("H2").Value = 168162599
$ws.Range("I2").Value = "TEAM"
$ws.Range("D3").Value = 21
$ws.Range("E3").Value = 30.07999992370605
$ws.Range("F3").Value = 31.45999908447266
$ws.Range("G3").Value = 21
$ws.Range("H3").Value = 168162599
$ws.Range("I3").Value = "TEAM"
$ws.Range("D4").Value = 21
$ws.Range("E4").Value = 30.07999992370605
$ws.Range("F4").Value = 31.45999908447266
$ws.Range("G4").Value = 21
$ws.Range("H4").Value = 168162599
$ws.Range("I4").Value = "TEAM"
$ws.Range("D5").Value = 28.94000053405762
$ws.Range("E5").Value = 20.77000045776367
$ws.Range("F5").Value = 29.03000068664551
$ws.Range("G5").Value = 20.17000007629395
$ws.Range("H5").Value = 168162599
$ws.Range("I5").Value = "TEAM"
$ws.Range("D6").Value = 24.85000038146973
$ws.Range("E6").Value = 23.17000007629395
$ws.Range("F6").Value = 24.97500038146973
$ws.Range("G6").Value = 22.5
$ws.Range("H6").Value = 168162599
$ws.Range("I6").Value = "TEAM"
$ws.Range("D7").Value = 25.68000030517578
$ws.Range("E7").Value = 29.96999931335449
$ws.Range("F7").Value = 30
$ws.Range("G7").Value = 25.21999931335449
$ws.Range("H7").Value = 168162599
$ws.Range("I7").Value = "TEAM"
$ws.Range("D8").Value = 29.79999923706055
$ws.Range("E8").Value = 26.86000061035156
$ws.Range("F8").Value = 30
$ws.Range("G8").Value = 26.05999946594238
$ws.Range("H8").Value = 168162599
$ws.Range("I8").Value = "TEAM"
$ws.Range("D9").Value = 24.31999969482422
$ws.Range("E9").Value = 27.6299991607666
$ws.Range("F9").Value = 28.69000053405762
$ws.Range("G9").Value = 24.20000076293945
$ws.Range("H9").Value = 168162599
$ws.Range("I9").Value = "TEAM"
$ws.Range("D10").Value = 29.85000038146973
$ws.Range("E10").Value = 34.47999954223633
$ws.Range("F10").Value = 34.97000122070312
$ws.Range("G10").Value = 29.55999946594238
$ws.Range("H10").Value = 168162599
$ws.Range("I10").Value = "TEAM"
$ws.Range("D11").Value = 35.25
$ws.Range("E11").Value = 35.81999969482422
$ws.Range("F11").Value = 39.25
$ws.Range("G11").Value = 33.34000015258789
$ws.Range("H11").Value = 168162599
$ws.Range("I11").Value = "TEAM"
$ws.Range("D12").Value = 35
$ws.Range("E12").Value = 48.36999893188477
$ws.Range("F12").Value = 50.88000106811523
$ws.Range("G12").Value = 35
$ws.Range("H12").Value = 168162599
$ws.Range("I12").Value = "TEAM"
$ws.Range("D13").Value = 45.93999862670898
$ws.Range("E13").Value = 53.9900016784668
$ws.Range("F13").Value = 55.90999984741211
$ws.Range("G13").Value = 45.61000061035156
$ws.Range("H13").Value = 168162599
$ws.Range("I13").Value = "TEAM"
$ws.Range("D14").Value = 53.54000091552734
$ws.Range("E14").Value = 55.97999954223633
$ws.Range("F14").Value = 62.77999877929688
$ws.Range("G14").Value = 51.7599983215332
$ws.Range("H14").Value = 168162599
$ws.Range("I14").Value = "TEAM"
$ws.Range("D15").Value = 62.02999877929688
$ws.Range("E15").Value = 72.41000366210938
$ws.Range("F15").Value = 79.81999969482422
$ws.Range("G15").Value = 61.65999984741211
$ws.Range("H15").Value = 168162599
$ws.Range("I15").Value = "TEAM"
$ws.Range("D16").Value = 97.18000030517578
$ws.Range("E16").Value = 75.91000366210938
$ws.Range("F16").Value = 98.21299743652344
$ws.Range("G16").Value = 66.80000305175781
$ws.Range("H16").Value = 168162599
$ws.Range("I16").Value = "TEAM"
$ws.Range("D17").Value = 85.81999969482422
$ws.Range("E17").Value = 98.40000152587891
$ws.Range("F17").Value = 100
$ws.Range("G17").Value = 84.18000030517578
$ws.Range("H17").Value = 168162599
$ws.Range("I17").Value = "TEAM"
$ws.Range("D18").Value = 114.0299987792969
$ws.Range("E18").Value = 110.1500015258789
$ws.Range("F18").Value = 117.0599975585938
$ws.Range("G18").Value = 100.25
$ws.Range("H18").Value = 168162599
$ws.Range("I18").Value = "TEAM"
$ws.Range("D19").Value = 133.1999969482422
$ws.Range("E19").Value = 140.1199951171875
$ws.Range("F19").Value = 149.8000030517578
$ws.Range("G19").Value = 129.6000061035156
$ws.Range("H19").Value = 168162599
$ws.Range("I19").Value = "TEAM"
$ws.Range("D20").Value = 125.2799987792969
$ws.Range("E20").Value = 120.7900009155273
$ws.Range("F20").Value = 133.6900024414062
$ws.Range("G20").Value = 107
$ws.Range("H20").Value = 168162599
$ws.Range("I20").Value = "TEAM"
$ws.Range("D21").Value = 121.1500015258789
$ws.Range("E21").Value = 147
$ws.Range("F21").Value = 151.8699951171875
$ws.Range("G21").Value = 119.8399963378906
$ws.Range("H21").Value = 168162599
$ws.Range("I21").Value = "TEAM"
$ws.Range("D22").Value = 132.9600067138672
$ws.Range("E22").Value = 155.4900054931641
$ws.Range("F22").Value = 158.97900390625
$ws.Range("G22").Value = 126.5400009155273
$ws.Range("H22").Value = 168162599
$ws.Range("I22").Value = "TEAM"
$ws.Range("D23").Value = 180.9900054931641
$ws.Range("E23").Value = 176.6499938964844
$ws.Range("F23").Value = 198.4100036621093
$ws.Range("G23").Value = 169.1100006103516
$ws.Range("H23").Value = 168162599
$ws.Range("I23").Value = "TEAM"
$ws.Range("D24").Value = 183.2200012207031
$ws.Range("E24").Value = 191.6199951171875
$ws.Range("F24").Value = 216.2949981689453
$ws.Range("G24").Value = 180.5700073242188
$ws.Range("H24").Value = 168162599
$ws.Range("I24").Value = "TEAM"
$ws.Range("D25").Value = 234.5
$ws.Range("E25").Value = 231.1300048828125
$ws.Range("F25").Value = 240.259994506836
$ws.Range("G25").Value = 208.6199951171875
$ws.Range("H25").Value = 168162599
$ws.Range("I25").Value = "TEAM"
$ws.Range("D26").Value = 216.009994506836
$ws.Range("E26").Value = 237.5599975585937
$ws.Range("F26").Value = 244.9299926757812
$ws.Range("G26").Value = 212.25
$ws.Range("H26").Value = 168162599
$ws.Range("I26").Value = "TEAM"
$ws.Range("D27").Value = 256.2300109863281
$ws.Range("E27").Value = 325.1199951171875
$ws.Range("F27").Value = 349.5
$ws.Range("G27").Value = 255.3500061035156
$ws.Range("H27").Value = 168162599
$ws.Range("I27").Value = "TEAM"
$ws.Range("D28").Value = 393
$ws.Range("E28").Value = 458.1300048828125
$ws.Range("F28").Value = 483.1300048828125
$ws.Range("G28").Value = 376
$ws.Range("H28").Value = 168162599
$ws.Range("I28").Value = "TEAM"
$ws.Range("D29").Value = 380.4299926757813
$ws.Range("E29").Value = 324.3399963378906
$ws.Range("F29").Value = 380.6300048828125
$ws.Range("G29").Value = 273.4219970703125
$ws.Range("H29").Value = 168162599
$ws.Range("I29").Value = "TEAM"
$ws.Range("D30").Value = 297.0299987792969
$ws.Range("E30").Value = 224.8300018310547
$ws.Range("F30").Value = 318.6400146484375
$ws.Range("G30").Value = 224.0299987792969
$ws.Range("H30").Value = 168162599
$ws.Range("I30").Value = "TEAM"
$ws.Range("D31").Value = 188.3899993896484
$ws.Range("E31").Value = 209.3200073242188
$ws.Range("F31").Value = 220.1300048828125
$ws.Range("G31").Value = 178.6799926757812
$ws.Range("H31").Value = 168162599
$ws.Range("I31").Value = "TEAM"
$ws.Range("D32").Value = 214.5500030517578
$ws.Range("E32").Value = 202.729995727539
$ws.Range("F32").Value = 246.1499938964844
$ws.Range("G32").Value = 184.4700012207031
$ws.Range("H32").Value = 168162599
$ws.Range("I32").Value = "TEAM"
$ws.Range("D33").Value = 131.9100036621094
$ws.Range("E33").Value = 161.6199951171875
$ws.Range("F33").Value = 165.0599975585938
$ws.Range("G33").Value = 116.4000015258789
$ws.Range("H33").Value = 168162599
$ws.Range("I33").Value = "TEAM"
$ws.Range("D34").Value = 169.3200073242188
$ws.Range("E34").Value = 147.6600036621094
$ws.Range("F34").Value = 170.6699981689453
$ws.Range("G34").Value = 142.8600006103516
$ws.Range("H34").Value = 168162599
$ws.Range("I34").Value = "TEAM"
$ws.Range("D35").Value = 168.1699981689453
$ws.Range("E35").Value = 181.9400024414062
$ws.Range("F35").Value = 195.9900054931641
$ws.Range("G35").Value = 162.5899963378906
$ws.Range("H35").Value = 168162599
$ws.Range("I35").Value = "TEAM"
$ws.Range("D36").Value = 201.8000030517578
$ws.Range("E36").Value = 180.6399993896484
$ws.Range("F36").Value = 208.8600006103516
$ws.Range("G36").Value = 173.5809936523438
$ws.Range("H36").Value = 168162599
$ws.Range("I36").Value = "TEAM"
$ws.Range("D37").Value = 235.009994506836
$ws.Range("E37").Value = 249.7700042724609
$ws.Range("F37").Value = 258.6900024414062
$ws.Range("G37").Value = 218.8500061035156
$ws.Range("H37").Value = 168162599
$ws.Range("I37").Value = "TEAM"
$ws.Range("D38").Value = 195.1999969482422
$ws.Range("E38").Value = 172.3000030517578
$ws.Range("F38").Value = 217.4600067138672
$ws.Range("G38").Value = 171
$ws.Range("H38").Value = 168162599
$ws.Range("I38").Value = "TEAM"
$ws.Range("D39").Value = 176.8800048828125
$ws.Range("E39").Value = 176.5700073242188
$ws.Range("F39").Value = 188.0599975585937
$ws.Range("G39").Value = 167.25
$ws.Range("H39").Value = 168162599
$ws.Range("I39").Value = "TEAM"
$ws.Range("D40").Value = 163.2299957275391
$ws.Range("E40").Value = 188.5399932861328
$ws.Range("F40").Value = 196.25
$ws.Range("G40").Value = 155.3699951171875
$ws.Range("H40").Value = 168162599
$ws.Range("I40").Value = "TEAM"
$ws.Range("D41").Value = 244.5299987792969
$ws.Range("E41").Value = 306.7799987792969
$ws.Range("F41").Value = 324.3699951171875
$ws.Range("G41").Value = 235.4400024414062
$ws.Range("H41").Value = 168162599
$ws.Range("I41").Value = "TEAM"
$ws.Range("D42").Value = 211.1699981689453
$ws.Range("E42").Value = 228.3099975585937
$ws.Range("F42").Value = 233.1300048828125
$ws.Range("G42").Value = 173.4589996337891
$ws.Range("H42").Value = 168162599
$ws.Range("I42").Value = "TEAM"
$ws.Range("D43").Value = 203.6600036621093
$ws.Range("E43").Value = 191.7799987792969
$ws.Range("F43").Value = 222.5899963378907
$ws.Range("G43").Value = 185.0800018310547
$ws.Range("H43").Value = 168162599
$ws.Range("I43").Value = "TEAM"
